$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# HOUR_APPR_PROCESS_START column is V, data rows 2-19.
# Convert the numeric hour values to text strings formatted as "H:00:00"
for ($row = 2; $row -le 19; $row++) {
    $cell = $ws.Cells.Item($row, 22)  # column V = 22
    $hour = $cell.Value2
    $cell.Value = [string]$hour + ":00:00"
}
